$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 17
$ws.Range("F3").Value = 977
$ws.Range("F4").Value = 222
$ws.Range("F6").Value = 1104
$ws.Range("F7").Value = 883
$ws.Range("F8").Value = 269
$ws.Range("F9").Value = 60
$ws.Range("F11").Value = 866
$ws.Range("F12").Value = 303
$ws.Range("F13").Value = 589
$ws.Range("F14").Value = 514
$ws.Range("F15").Value = 1356
$ws.Range("G15").Value = 60
$ws.Range("F17").Value = 1278
$ws.Range("F18").Value = 1231
$ws.Range("F19").Value = 2911
$ws.Range("F20").Value = 1496
$ws.Range("F21").Value = 734
$ws.Range("F22").Value = 207
$ws.Range("F23").Value = 1295
$ws.Range("F25").Value = 1050
$ws.Range("F26").Value = 367
$ws.Range("F27").Value = 3224
$ws.Range("F28").Value = 631
$ws.Range("F29").Value = 543
$ws.Range("F30").Value = 1439

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 42
$ws.Range("F8").Value = 31
$ws.Range("F12").Value = 8

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 17
$ws.Range("F6").Value = 977
$ws.Range("F7").Value = 222
$ws.Range("F10").Value = 1104
$ws.Range("F11").Value = 883
$ws.Range("F12").Value = 269
$ws.Range("F14").Value = 60
$ws.Range("F15").Value = 42
$ws.Range("F16").Value = 42
$ws.Range("F18").Value = 31
$ws.Range("F21").Value = 866
$ws.Range("F22").Value = 303
$ws.Range("F23").Value = 589
$ws.Range("F24").Value = 514
$ws.Range("F25").Value = 1356
$ws.Range("G25").Value = 60
$ws.Range("F27").Value = 1278
$ws.Range("F28").Value = 1231
$ws.Range("F29").Value = 2911
$ws.Range("F30").Value = 1496
$ws.Range("F31").Value = 734
$ws.Range("F32").Value = 207
$ws.Range("F33").Value = 1295
$ws.Range("F37").Value = 1050
$ws.Range("F38").Value = 367
$ws.Range("F39").Value = 3224
$ws.Range("F40").Value = 631
$ws.Range("F41").Value = 543
$ws.Range("F42").Value = 1439
$ws.Range("F43").Value = 8

Write-Host "Applied 59 cell updates across 展览, 演出, 全部类型 sheets"
